$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("6x6")
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
